$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 2 (existing rows 2..258 shift down to 6..262)
$ws.Rows.Item(2).Resize(4).Insert()

# The inserted rows inherit formatting from the header row above; strip that so
# the new rows look like the plain data rows that follow them.
$newRowsRange = $ws.Range("A2:F5")
$newRowsRange.ClearFormats()
$newRowsRange.ClearContents()
$ws.Range("A2:A5,D2:D5,E2:E5,F2:F5").Clear()

# Fill in the 4 newly inserted rows with the new songs
$ws.Cells.Item(2, 2).Value = "Taylor Swift - Cruel Summer (Lyrics)"
$ws.Cells.Item(2, 3).Value = "https://www.youtube.com/watch?v=P8T1rUpVdXE"

$ws.Cells.Item(3, 2).Value = "Alan Walker - Lily (Lyrics) ft. K391, Emelie Hollow"
$ws.Cells.Item(3, 3).Value = "https://www.youtube.com/watch?v=sZ5OUc7Ccwo"

$ws.Cells.Item(4, 2).Value = "Sebastián Yatra - Adiós (Letra)"
$ws.Cells.Item(4, 3).Value = "https://www.youtube.com/watch?v=DULoaFTPB60"

$ws.Cells.Item(5, 2).Value = "Taylor Swift - Out Of The Woods (Lyrics)"
$ws.Cells.Item(5, 3).Value = "https://www.youtube.com/watch?v=FEzj8K1cL6w"
